$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain decimal number must be forced to
# Text format first, otherwise Excel auto-converts the literal (e.g.
# "587.80") into a numeric value and the exact text (with trailing
# zeros) is lost.
$forceTextCells = @("D5", "D6", "D8", "D9", "D10", "D11", "D13", "D14", "D18", "D19", "D21", "D23", "D24", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D34", "D35", "D38", "D39", "D40", "D41", "D45", "D50")
foreach ($cellRef in $forceTextCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range('D2').Value = '69.357.76'
$ws.Range('E2').Value = '  +2.06%  '
$ws.Range('D3').Value = '3.387.33'
$ws.Range('E3').Value = '  +1.40%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '587.80'
$ws.Range('D6').Value = '179.71'
$ws.Range('E6').Value = '  +1.08%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').Value = '0.597'
$ws.Range('E8').Value = '  +1.19%  '
$ws.Range('D9').Value = '0.194'
$ws.Range('E9').Value = '  +4.88%  '
$ws.Range('D10').Value = '0.591'
$ws.Range('E10').Value = '  +1.24%  '
$ws.Range('D11').Value = '48.44'
$ws.Range('E11').Value = '  +2.45%  '
$ws.Range('E12').Value = '  +2.80%  '
$ws.Range('D13').Value = '678.82'
$ws.Range('E13').Value = '  -3.98%  '
$ws.Range('D14').Value = '8.61'
$ws.Range('D15').Value = '3.929.69'
$ws.Range('E15').Value = '  +1.23%  '
$ws.Range('D16').Value = '69.421.83'
$ws.Range('E16').Value = '  +2.10%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '3.392.80'
$ws.Range('E17').Value = '  +1.36%  '
$ws.Range('B18').Value = 'TRON'
$ws.Range('C18').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D18').Value = '0.120'
$ws.Range('E18').Value = '  +1.62%  '
$ws.Range('D19').Value = '17.67'
$ws.Range('E19').Value = '  +0.61%  '
$ws.Range('E20').Value = '  +1.68%  '
$ws.Range('D21').Value = '0.905'
$ws.Range('E21').Value = '  +0.82%  '
$ws.Range('E22').Value = '  +0.32%  '
$ws.Range('D23').Value = '17.23'
$ws.Range('E23').Value = '  +0.70%  '
$ws.Range('D24').Value = '103.58'
$ws.Range('E24').Value = '  +2.98%  '
$ws.Range('E25').Value = '  +0.11%  '
$ws.Range('D26').Value = '2.72'
$ws.Range('E26').Value = '  +0.48%  '
$ws.Range('D27').Value = '9.68'
$ws.Range('E27').Value = '  +0.37%  '
$ws.Range('D28').Value = '34.06'
$ws.Range('E28').Value = '  +2.64%  '
$ws.Range('D29').Value = '8.73'
$ws.Range('E29').Value = '  +1.44%  '
$ws.Range('D30').Value = '6.96'
$ws.Range('E30').Value = '  -1.16%  '
$ws.Range('B31').Value = 'Bittensor'
$ws.Range('C31').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D31').Value = '561.27'
$ws.Range('E31').Value = '  -1.60%  '
$ws.Range('B32').Value = 'Cosmos'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D32').Value = '11.14'
$ws.Range('E32').Value = '  +0.99%  '
$ws.Range('E33').Value = '  +0.35%  '
$ws.Range('D34').Value = '3.57'
$ws.Range('E34').Value = '  +4.56%  '
$ws.Range('D35').Value = '58.56'
$ws.Range('E35').Value = '  +1.27%  '
$ws.Range('E36').Value = '  +0.17%  '
$ws.Range('D37').Value = '3.688.07'
$ws.Range('E37').Value = '  -0.30%  '
$ws.Range('B38').Value = 'InjectiveProtocol'
$ws.Range('C38').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D38').Value = '35.59'
$ws.Range('E38').Value = '  +2.22%  '
$ws.Range('B39').Value = 'Kaspa'
$ws.Range('C39').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D39').Value = '0.139'
$ws.Range('E39').Value = '  +4.40%  '
$ws.Range('D40').Value = '3.26'
$ws.Range('E40').Value = '  +2.43%  '
$ws.Range('D41').Value = '2.69'
$ws.Range('E41').Value = '  +1.15%  '
$ws.Range('D42').Value = '0.0₃0699'
$ws.Range('E42').Value = '  +2.85%  '
$ws.Range('E43').Value = '  +0.59%  '
$ws.Range('E44').Value = '  +3.61%  '
$ws.Range('D45').Value = '3.27'
$ws.Range('E45').Value = '  -2.20%  '
$ws.Range('E46').Value = '  -0.13%  '
$ws.Range('E47').Value = '  +0.90%  '
$ws.Range('E48').Value = '  +4.76%  '
$ws.Range('E49').Value = '  +0.00%  '
$ws.Range('D50').Value = '133.30'
$ws.Range('E50').Value = '  +1.68%  '
$ws.Range('E51').Value = '  +3.15%  '
